$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.660.36'
$ws.Range('E2').Value = '  +1.46%  '
$ws.Range('E3').Value = '  -0.65%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '591.42'
$ws.Range('E5').Value = '  -0.14%  '
$ws.Range('D6').Value = '136.56'
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '3.177.28'
$ws.Range('E8').Value = '  -0.78%  '
$ws.Range('D9').Value = '0.515'
$ws.Range('E9').Value = '  +1.58%  '
$ws.Range('E10').Value = '  -1.09%  '
$ws.Range('E11').Value = '  +0.51%  '
$ws.Range('D12').Value = '0.457'
$ws.Range('E12').Value = '  +0.28%  '
$ws.Range('E13').Value = '  +0.78%  '
$ws.Range('D14').Value = '34.89'
$ws.Range('E14').Value = '  +3.92%  '
$ws.Range('D15').Value = '3.702.62'
$ws.Range('E15').Value = '  -0.74%  '
$ws.Range('E16').Value = '  -0.55%  '
$ws.Range('D17').Value = '3.180.97'
$ws.Range('E17').Value = '  -0.63%  '
$ws.Range('D18').Value = '63.650.17'
$ws.Range('E18').Value = '  +1.34%  '
$ws.Range('D19').Value = '6.58'
$ws.Range('E19').Value = '  -2.00%  '
$ws.Range('D20').Value = '463.26'
$ws.Range('E20').Value = '  -0.32%  '
$ws.Range('D21').Value = '13.99'
$ws.Range('E21').Value = '  +0.60%  '
$ws.Range('D22').Value = '0.701'
$ws.Range('D23').Value = '7.69'
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').Value = '13.29'
$ws.Range('E24').Value = '  -0.99%  '
$ws.Range('D25').Value = '83.41'
$ws.Range('E25').Value = '  -0.70%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('E27').Value = '  -0.83%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').Value = '6.88'
$ws.Range('E29').Value = '  -0.43%  '
$ws.Range('D30').Value = '2.09'
$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('D31').Value = '7.80'
$ws.Range('E31').Value = '  -1.17%  '
$ws.Range('E32').Value = '  -0.37%  '
$ws.Range('D33').Value = '0.102'
$ws.Range('E33').Value = '  -0.75%  '
$ws.Range('D34').Value = '2.43'
$ws.Range('E34').Value = '  -0.16%  '
$ws.Range('E35').Value = '  -1.58%  '
$ws.Range('D36').Value = '5.91'
$ws.Range('E36').Value = '  +0.81%  '
$ws.Range('D37').Value = '0.0₃0742'
$ws.Range('E37').Value = '  +5.83%  '
$ws.Range('D38').Value = '51.65'
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('E39').Value = '  +0.78%  '
$ws.Range('D40').Value = '8.16'
$ws.Range('E40').Value = '  +0.65%  '
$ws.Range('E41').Value = '  +2.01%  '
$ws.Range('E42').Value = '  -1.04%  '
$ws.Range('D43').Value = '399.44'
$ws.Range('E43').Value = '  -4.93%  '
$ws.Range('D44').Value = '2.792.23'
$ws.Range('E44').Value = '  -7.01%  '
$ws.Range('D45').Value = '0.254'
$ws.Range('E45').Value = '  -0.43%  '
$ws.Range('E46').Value = '  -0.17%  '
$ws.Range('D47').Value = '127.66'
$ws.Range('E47').Value = '  +1.93%  '
$ws.Range('B48').Value = 'Arweave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D48').Value = '35.89'
$ws.Range('E48').Value = '  +1.34%  '
$ws.Range('B49').Value = 'USDe'
$ws.Range('C49').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D49').Value = '0.999'
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('D50').Value = '25.62'
$ws.Range('E50').Value = '  -0.23%  '
$ws.Range('E51').Value = '  -0.12%  '
